$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 71: convert E71, F71, I71 from text to real numbers ---
$ws.Range("E71").Value = 213
$ws.Range("F71").Value = 123
$ws.Range("I71").Value = 345

# --- Row 72: new row of (mostly textual) data ---
# Use a leading apostrophe to force Excel to store these number-looking /
# date-looking values as literal text instead of auto-converting them to
# numbers or dates, then clear the resulting "quote prefix" formatting so
# no extra cell style is left behind on the cell.
$ws.Range("A72").Value = "'2024-04-30"
$ws.Range("A72").ClearFormats()

$ws.Range("B72").Value = "Flowserve"

$ws.Range("C72").Value = "'234"
$ws.Range("C72").ClearFormats()

$ws.Range("D72").Value = "'234"
$ws.Range("D72").ClearFormats()

$ws.Range("E72").Value = "'234"
$ws.Range("E72").ClearFormats()

$ws.Range("F72").Value = "'234"
$ws.Range("F72").ClearFormats()

$ws.Range("G72").Value = "'234"
$ws.Range("G72").ClearFormats()

$ws.Range("H72").Value = "Acessos"

$ws.Range("I72").Value = "'234"
$ws.Range("I72").ClearFormats()

# J72 is left blank (matches the empty inline string cell in the target)
